$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change cell A2 value from "two" to "one" (adds a new shared string)
$ws.Range("A2").Value = "one"

# Update the selection to A3 (matches the saved sheetView selection)
$ws.Range("A3").Select()
